$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that need to be written as literal TEXT (not auto-converted to numbers)
# are first formatted as Text ("@"), written, then restored to the default "Normal"
# style so no stray number-format / quote-prefix styling is left behind.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "24.626.16"
$ws.Range("E2").Value = "  -0.37%  "

# Row 3
Set-TextValue "D3" "1.676.09"
$ws.Range("E3").Value = "  -0.56%  "

# Row 4
Set-TextValue "D4" "1.006"
$ws.Range("E4").Value = "  +0.46%  "

# Row 5
Set-TextValue "D5" "307.33"
$ws.Range("E5").Value = "  +0.37%  "

# Row 6
Set-TextValue "D6" "0.9998"
$ws.Range("E6").Value = "  +0.40%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
Set-TextValue "D8" "48.22"
$ws.Range("E8").Value = "  -2.66%  "

# Row 9
Set-TextValue "D9" "0.3367"
$ws.Range("E9").Value = "  -1.69%  "

# Row 10
Set-TextValue "D10" "1.176"
$ws.Range("E10").Value = "  +1.00%  "

# Row 11
Set-TextValue "D11" "0.07322"
$ws.Range("E11").Value = "  +1.18%  "

# Row 12
Set-TextValue "D12" "1.001"
$ws.Range("E12").Value = "  +0.37%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D13" "6.176"
$ws.Range("E13").Value = "  +1.27%  "

# Row 14
$ws.Range("B14").Value = "Solana"
$ws.Range("C14").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D14" "20.49"
$ws.Range("E14").Value = "  +1.56%  "

# Row 15
Set-TextValue "D15" "6.793"
$ws.Range("E15").Value = "  +1.41%  "

# Row 16
Set-TextValue "D16" "1.676.48"
$ws.Range("E16").Value = "  -0.32%  "

# Row 17
Set-TextValue "D17" "0.00001097"
$ws.Range("E17").Value = "  -0.52%  "

# Row 18
Set-TextValue "D18" "0.06632"
$ws.Range("E18").Value = "  -0.48%  "

# Row 19
Set-TextValue "D19" "0.9998"
$ws.Range("E19").Value = "  +0.45%  "

# Row 20
Set-TextValue "D20" "81.66"
$ws.Range("E20").Value = "  +0.87%  "

# Row 21
Set-TextValue "D21" "16.80"
$ws.Range("E21").Value = "  +2.62%  "

# Row 22
Set-TextValue "D22" "6.195"
$ws.Range("E22").Value = "  +1.94%  "

# Row 23
Set-TextValue "D23" "12.67"
$ws.Range("E23").Value = "  +5.01%  "

# Row 24
Set-TextValue "D24" "24.632.93"
$ws.Range("E24").Value = "  +0.00%  "

# Row 25
Set-TextValue "D25" "2.435"
$ws.Range("E25").Value = "  +0.83%  "

# Row 26
Set-TextValue "D26" "2.694"
$ws.Range("E26").Value = "  +1.15%  "

# Row 27
Set-TextValue "D27" "19.83"
$ws.Range("E27").Value = "  +1.99%  "

# Row 28
Set-TextValue "D28" "149.29"
$ws.Range("E28").Value = "  -2.05%  "

# Row 29
Set-TextValue "D29" "130.19"
$ws.Range("E29").Value = "  +2.04%  "

# Row 30
Set-TextValue "D30" "1.860.07"
$ws.Range("E30").Value = "  -0.44%  "

# Row 31
Set-TextValue "D31" "1.215"
$ws.Range("E31").Value = "  +24.06%  "

# Row 32
Set-TextValue "D32" "6.500"
$ws.Range("E32").Value = "  +3.87%  "

# Row 33
Set-TextValue "D33" "4.158"
$ws.Range("E33").Value = "  +3.22%  "

# Row 34
Set-TextValue "D34" "0.08592"
$ws.Range("E34").Value = "  +2.17%  "

# Row 35
Set-TextValue "D35" "13.32"
$ws.Range("E35").Value = "  +7.87%  "

# Row 36
Set-TextValue "D36" "1.725"
$ws.Range("E36").Value = "  +2.34%  "

# Row 37
Set-TextValue "D37" "5.413"
$ws.Range("E37").Value = "  +2.09%  "

# Row 38
Set-TextValue "D38" "0.06453"
$ws.Range("E38").Value = "  +1.65%  "

# Row 39
Set-TextValue "D39" "8.849"
$ws.Range("E39").Value = "  +2.43%  "

# Row 40
Set-TextValue "D40" "0.02351"
$ws.Range("E40").Value = "  +1.85%  "

# Row 41
Set-TextValue "D41" "0.2163"
$ws.Range("E41").Value = "  +3.76%  "

# Row 42
Set-TextValue "D42" "1.237"
$ws.Range("E42").Value = "  -0.59%  "

# Row 43
Set-TextValue "D43" "0.6250"
$ws.Range("E43").Value = "  +2.67%  "

# Row 44
$ws.Range("E44").Value = "  +0.46%  "

# Row 45
Set-TextValue "D45" "13.42"
$ws.Range("E45").Value = "  +2.54%  "

# Row 46
Set-TextValue "D46" "3.778"
$ws.Range("E46").Value = "  +0.74%  "

# Row 47
Set-TextValue "D47" "0.5946"
$ws.Range("E47").Value = "  +1.19%  "

# Row 48
$ws.Range("E48").Value = "  +2.34%  "

# Row 49
Set-TextValue "D49" "125.77"
$ws.Range("E49").Value = "  +0.20%  "

# Row 50
Set-TextValue "D50" "0.07130"
$ws.Range("E50").Value = "  -1.47%  "

# Row 51
Set-TextValue "D51" "76.99"
$ws.Range("E51").Value = "  +1.64%  "
